$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 11 A:E fully (values + formatting), then set A11 content
$ws.Range("A11:E11").Clear()
$ws.Range("A11").Value = "No good drivers found."

# Delete rows 12-16
$ws.Range("A12:J16").Delete()

# set column widths
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(5).ColumnWidth = 2
